$d = $word.ActiveDocument

$replacements = @(
    @{ old = "2025-08-29 Friday"; new = "2025-08-30 Saturday" },
    @{ old = "113×2="; new = "268×5=" },
    @{ old = "875×9="; new = "794×9=" },
    @{ old = "399×4="; new = "621×7=" },
    @{ old = "459×3="; new = "250×6=" },
    @{ old = "432×3="; new = "679×2=" },
    @{ old = "146×2="; new = "624×8=" },
    @{ old = "972×5="; new = "529×9=" },
    @{ old = "941×3="; new = "466×5=" },
    @{ old = "372×3="; new = "336×8=" },
    @{ old = "683×3="; new = "536×4=" },
    @{ old = "820×6="; new = "291×3=" },
    @{ old = "703×7="; new = "613×2=" },
    @{ old = "997×8="; new = "257×5=" },
    @{ old = "408×2="; new = "593×4=" },
    @{ old = "958×4="; new = "597×4=" },
    @{ old = "423×9="; new = "275×5=" },
    @{ old = "209×8="; new = "112×6=" },
    @{ old = "422×6="; new = "856×9=" },
    @{ old = "991×5="; new = "456×6=" },
    @{ old = "394×3="; new = "439×5=" },
    @{ old = "986×9="; new = "922×2=" },
    @{ old = "348×3="; new = "354×4=" },
    @{ old = "403×7="; new = "736×6=" },
    @{ old = "300×2="; new = "582×2=" },
    @{ old = "732×7="; new = "150×4=" }
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
